# Updated car model drag coeff and added low drag configs

$wb = $excel.ActiveWorkbook

$infoSheet = $wb.Worksheets.Item("Info")

# --- Info sheet value/formula updates ---

# Lift Coefficient CL: -2 -> -1.98
$infoSheet.Range("C8").Value = -1.98

# Drag Coefficient CD: -1.2 -> -1.33
$infoSheet.Range("C9").Value = -1.33

# Front Aero Distribution: 47 -> formula =100-56.3 (43.7)
$infoSheet.Range("C12").Formula = "=100-56.3"

# Frontal Area: 1.1 -> 1.15
$infoSheet.Range("C13").Value = 1.15

# --- Selection / active sheet updates ---

# Make "Info" the active / selected tab (tabSelected moves from "Torque Curve" to "Info")
$infoSheet.Activate()

# Select C12:C13 on Info sheet, with active cell C12
$infoSheet.Range("C12:C13").Select()
